# Add new Parse Examples
# Target worksheet: "Игрушки" (2nd tab), column F holds "Parse Examples" annotations
# next to the URL column D.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Игрушки")

# Row 22: new "Watermark" example
$ws.Range("F22").Value = "Watermark"

# Rows 31-43: new "No price" examples
$ws.Range("F31").Value = "No price"
$ws.Range("F32").Value = "No price"
$ws.Range("F33").Value = "No price"
$ws.Range("F34").Value = "No price"
$ws.Range("F35").Value = "No price"
$ws.Range("F36").Value = "No price"
$ws.Range("F37").Value = "No price"
$ws.Range("F38").Value = "No price"
$ws.Range("F39").Value = "No price"
$ws.Range("F40").Value = "No price"
$ws.Range("F41").Value = "No price"
$ws.Range("F42").Value = "No price"
$ws.Range("F43").Value = "No price"

# Update the view state to match the author's final scroll/selection position.
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 15
$ws.Range("H40").Select()
